$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 949
$ws.Range("I11").Value = 949
$ws.Range("K11").Value = 949
$ws.Range("M11").Value = -809
$ws.Range("H31").Value = 2037
$ws.Range("I31").Value = 2037
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6111
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = -5881
$ws.Range("H38").Value = 518.6
$ws.Range("I38").Value = 596.9231
$ws.Range("J38").Value = 9.5
$ws.Range("K38").Value = 1790.7693
$ws.Range("L38").Value = 28.5
$ws.Range("M38").Value = -1418.7693
$ws.Range("N38").Value = -772.5
$ws.Range("H53").Value = 409.77777
$ws.Range("I53").Value = 57.4
$ws.Range("K53").Value = 57.4
$ws.Range("M53").Value = 579.6
$ws.Range("H64").Value = 4626.273
$ws.Range("I64").Value = 4484.143
$ws.Range("K64").Value = 4484.143
$ws.Range("M64").Value = -4236.143
$ws.Range("H67").Value = 4626.273
$ws.Range("I67").Value = 4484.143
$ws.Range("K67").Value = 4484.143
$ws.Range("M67").Value = -3626.143
$ws.Range("H74").Value = 45459670
$ws.Range("I74").Value = 50005240
$ws.Range("K74").Value = 50005240
$ws.Range("M74").Value = -50004304
$ws.Range("H77").Value = 45459670
$ws.Range("I77").Value = 50005240
$ws.Range("K77").Value = 250026200
$ws.Range("M77").Value = -250021520
$ws.Range("H118").Value = 1242
$ws.Range("I118").Value = 1242
$ws.Range("K118").Value = 3726
$ws.Range("M118").Value = -2069
$ws.Range("H131").Value = 775.3
$ws.Range("J131").Value = 498
$ws.Range("L131").Value = 1494
$ws.Range("N131").Value = -11574
$ws.Range("H137").Value = 2144.5
$ws.Range("I137").Value = 1501.8572
$ws.Range("J137").Value = 2490.5386
$ws.Range("K137").Value = 4505.571599999999
$ws.Range("L137").Value = 7471.6158
$ws.Range("M137").Value = -1955.571599999999
$ws.Range("N137").Value = -12571.6158
$ws.Range("H138").Value = 3180.0454
$ws.Range("I138").Value = 3095.0833
$ws.Range("K138").Value = 9285.249899999999
$ws.Range("M138").Value = -4145.249899999999
$ws.Range("H141").Value = 2571.8333
$ws.Range("I141").Value = 2571.8333
$ws.Range("K141").Value = 7715.499899999999
$ws.Range("M141").Value = -2535.499899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6060.1113
$ws.Range("I32").Value = 4144.96
$ws.Range("K32").Value = 4144.96
$ws.Range("M32").Value = -3857.96
$ws.Range("H61").Value = 43479564
$ws.Range("I61").Value = 47620190
$ws.Range("J61").Value = 2999.5
$ws.Range("K61").Value = 47620190
$ws.Range("L61").Value = 2999.5
$ws.Range("M61").Value = -47619978
$ws.Range("N61").Value = -3423.5
$ws.Range("H74").Value = 200006300
$ws.Range("I74").Value = 250005140
$ws.Range("K74").Value = 250005140
$ws.Range("M74").Value = -250004266
$ws.Range("H77").Value = 200006300
$ws.Range("I77").Value = 250005140
$ws.Range("K77").Value = 1250025700
$ws.Range("M77").Value = -1250021332
$ws.Range("H132").Value = 7146689.5
$ws.Range("I132").Value = 9093241
$ws.Range("K132").Value = 27279723
$ws.Range("M132").Value = -27277193
$ws.Range("H136").Value = 43479564
$ws.Range("I136").Value = 47620190
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 142860570
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = -142858020
$ws.Range("N136").Value = -14098.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1139.2858
$ws.Range("I20").Value = 995
$ws.Range("K20").Value = 995
$ws.Range("M20").Value = -748
$ws.Range("H24").Value = 7508
$ws.Range("I24").Value = 7508
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 7508
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = ""
$ws.Range("N24").Value = -7273
$ws.Range("H86").Value = 1912.1666
$ws.Range("I86").Value = 1990.9333
$ws.Range("K86").Value = 1990.9333
$ws.Range("M86").Value = -867.9332999999999
$ws.Range("H89").Value = 1912.1666
$ws.Range("I89").Value = 1990.9333
$ws.Range("K89").Value = 9954.666499999999
$ws.Range("M89").Value = -4338.666499999999
$ws.Range("H134").Value = 17862538
$ws.Range("I134").Value = 20005444
$ws.Range("K134").Value = 60016332
$ws.Range("M134").Value = -60013797

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 21745212
$ws.Range("I58").Value = 50011250
$ws.Range("J58").Value = 2105.7693
$ws.Range("K58").Value = 50011250
$ws.Range("L58").Value = 2105.7693
$ws.Range("M58").Value = -50011047
$ws.Range("N58").Value = -2511.7693
$ws.Range("H86").Value = 13181.5
$ws.Range("J86").Value = 15477.444
$ws.Range("L86").Value = 15477.444
$ws.Range("N86").Value = -17723.444
$ws.Range("H89").Value = 13181.5
$ws.Range("J89").Value = 15477.444
$ws.Range("L89").Value = 77387.22
$ws.Range("N89").Value = -88619.22
$ws.Range("H136").Value = 21745212
$ws.Range("I136").Value = 50011250
$ws.Range("J136").Value = 2105.7693
$ws.Range("K136").Value = 150033750
$ws.Range("L136").Value = 6317.3079
$ws.Range("M136").Value = -150031200
$ws.Range("N136").Value = -11417.3079

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 650.3077
$ws.Range("J33").Value = 801
$ws.Range("L33").Value = 4806
$ws.Range("N33").Value = -5372

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 11013
$ws.Range("I57").Value = 11013
$ws.Range("K57").Value = 11013
$ws.Range("M57").Value = -10193
$ws.Range("H70").Value = 4472.5
$ws.Range("I70").Value = 4168.857
$ws.Range("J70").Value = 4897.6
$ws.Range("K70").Value = 4168.857
$ws.Range("L70").Value = 4897.6
$ws.Range("M70").Value = -3898.857
$ws.Range("N70").Value = -5437.6
$ws.Range("H73").Value = 4472.5
$ws.Range("I73").Value = 4168.857
$ws.Range("J73").Value = 4897.6
$ws.Range("K73").Value = 4168.857
$ws.Range("L73").Value = 4897.6
$ws.Range("M73").Value = -3232.857
$ws.Range("N73").Value = -6769.6
$ws.Range("H113").Value = 65623.234
$ws.Range("I113").Value = 73706.336
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 73706.336
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -71536.336
$ws.Range("N113").Value = -9340

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2766.4443
$ws.Range("I22").Value = 3649.6667
$ws.Range("K22").Value = 3649.6667
$ws.Range("M22").Value = -3354.6667
$ws.Range("H27").Value = 2766.4443
$ws.Range("I27").Value = 3649.6667
$ws.Range("K27").Value = 3649.6667
$ws.Range("M27").Value = -3542.6667
$ws.Range("H40").Value = 3362.625
$ws.Range("I40").Value = 3100.3333
$ws.Range("K40").Value = 3100.3333
$ws.Range("M40").Value = -2964.3333
$ws.Range("H46").Value = 1174.8182
$ws.Range("J46").Value = 448.6
$ws.Range("L46").Value = 448.6
$ws.Range("N46").Value = -824.6
$ws.Range("H55").Value = 607.2143
$ws.Range("J55").Value = 1074.5
$ws.Range("L55").Value = 1074.5
$ws.Range("N55").Value = -1420.5
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").Value = ""
$ws.Range("H132").Value = 16554234
$ws.Range("I132").Value = 18463924
$ws.Range("K132").Value = 55391772
$ws.Range("M132").Value = -55389242
$ws.Range("H136").Value = 1381.8182
$ws.Range("I136").Value = 1245.3889
$ws.Range("J136").Value = 1995.75
$ws.Range("K136").Value = 3736.1667
$ws.Range("L136").Value = 5987.25
$ws.Range("M136").Value = -1186.1667
$ws.Range("N136").Value = -11087.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1056.2273
$ws.Range("I81").Value = 868.4286
$ws.Range("K81").Value = 1736.8572
$ws.Range("M81").Value = -675.8571999999999
$ws.Range("H84").Value = 1056.2273
$ws.Range("I84").Value = 868.4286
$ws.Range("K84").Value = 8684.286
$ws.Range("M84").Value = -3380.286
$ws.Range("H96").Value = 1505.5
$ws.Range("I96").Value = 1206.7
$ws.Range("K96").Value = 1206.7
$ws.Range("M96").Value = 166.3
$ws.Range("H104").Value = 41789.668
$ws.Range("J104").Value = 41789.668
$ws.Range("L104").Value = 41789.668
$ws.Range("N104").Value = -48777.668
$ws.Range("H107").Value = 709.8
$ws.Range("I107").Value = 724.75
$ws.Range("J107").Value = 650
$ws.Range("K107").Value = 2174.25
$ws.Range("L107").Value = 1950
$ws.Range("M107").Value = -254.25
$ws.Range("N107").Value = -5790
$ws.Range("H126").Value = 1632.8334
$ws.Range("I126").Value = 1862.25
$ws.Range("J126").Value = 1174
$ws.Range("K126").Value = 5586.75
$ws.Range("L126").Value = 3522
$ws.Range("M126").Value = -3116.75
$ws.Range("N126").Value = -8462
$ws.Range("H132").Value = 20009112
$ws.Range("I132").Value = 27783878
$ws.Range("K132").Value = 83351634
$ws.Range("M132").Value = -83349104
$ws.Range("H136").Value = 12821842
$ws.Range("I136").Value = 12821842
$ws.Range("K136").Value = 38465526
$ws.Range("M136").Value = -38462976
